$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WES")
$ws.Rows(2).Delete()
$ws.Range("B14:F14").ClearContents()
